$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the surviving UPN value (AmberR) to the new tenant id
$ws.Range("A2").Value = "AmberR@M365CPI14039056.OnMicrosoft.com"

# Remove the other two UPN rows (BillieV, admin) entirely
$ws.Range("A3:A4").Clear()

# Drop all hyperlinks on the sheet (A2/A3/A4 previously had mailto: links)
$ws.Hyperlinks.Delete()

# A2 no longer carries the Hyperlink style - reset it back to Normal
$ws.Range("A2").Style = "Normal"

# Move the active selection to A3 (matches the saved sheet view)
$ws.Range("A3").Select()
